# Insert a new daily price record for "Apio" (row 127) into the Vega Modelo
# de Temuco sheet, shifting all subsequent rows down by one (A1:R246 -> A1:R247).
#
# The newly inserted row is a copy of the row that used to be row 127
# (same market/region/quality/unit/origin/min-max-avg price & price-per-kg),
# except for the date (column D) and the volume (column J), which carry the
# new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 127

# Push rows 127:246 down to 128:247
$ws.Rows.Item($insertRow).Insert()

# The row that used to be at $insertRow now lives one row below; duplicate
# its contents into the freshly inserted (currently blank) row.
$ws.Range("A" + ($insertRow + 1) + ":R" + ($insertRow + 1)).Copy()
$ws.Range("A" + $insertRow).PasteSpecial()

# Apply the values that differ for the new record: date and volume.
$ws.Cells.Item($insertRow, 4).Value = 44601
$ws.Cells.Item($insertRow, 10).Value = 60
